# Sprint 40 - Day 2 Test Case Summary: fill in the counts for
# "Total testcase Written" / "Total Execution" / "Total Review"
# (email verification API test cases).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value  = 2000
$ws.Range("C10").Value = 1294
$ws.Range("C11").Value = 755

# Leave the cursor on the cell that was last edited instead of the
# previous scrolled-down selection.
$ws.Range("C11").Select()
